# Update cryptocurrency price/volume data on the active worksheet.
# Values are forced as literal text (matching the source inline-string cells)
# by using a leading apostrophe (Excel's quote-prefix convention), then the
# cell style is reset to Normal so no stray number-format/quote-prefix residue
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'68.435.38"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Formula = "'2.650.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +0.47%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Formula = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.15%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Formula = "'596.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -0.07%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Formula = "'159.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +3.02%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Formula = "'  +0.03%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Formula = "'0.541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'  -0.53%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Formula = "'  -0.10%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Formula = "'  -0.75%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Formula = "'5.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +0.76%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Formula = "'  +0.41%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Formula = "'28.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'  +1.33%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Formula = "'3.133.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +0.51%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Formula = "'  -1.61%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Formula = "'68.365.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +0.32%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Formula = "'2.639.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +0.17%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Formula = "'11.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +2.95%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Formula = "'364.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +0.53%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Formula = "'7.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +1.24%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Formula = "'  +1.86%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Formula = "'4.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -0.73%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Formula = "'2.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  +2.62%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Formula = "'75.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +0.05%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Formula = "'  +0.14%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Formula = "'9.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  +2.28%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E28").Formula = "'  -1.55%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Formula = "'  +0.12%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Formula = "'573.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +2.16%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Formula = "'  +0.92%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Formula = "'1.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +0.85%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Formula = "'  +1.02%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Formula = "'  +4.12%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Formula = "'  -0.44%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Formula = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  -0.05%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Formula = "'160.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +0.01%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Formula = "'19.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  +2.17%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Formula = "'0.370"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  -0.45%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Formula = "'  +0.27%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Formula = "'  +0.60%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Formula = "'  +0.70%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Formula = "'  -5.51%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Formula = "'  +0.13%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Formula = "'158.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  +1.03%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Formula = "'  +2.45%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Formula = "'21.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +0.94%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Formula = "'1.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  +1.06%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Formula = "'0.0780"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  -0.66%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Formula = "'0.576"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +2.90%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Formula = "'  +0.13%  "
$ws.Range("E51").Style = "Normal"
